$wb = $excel.ActiveWorkbook

$latestUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/444133beedd7e169e1d11d60f530615a0b07b21d/e2e/6d9ee57a-572e-4ce2-b42d-a8b49321472b.md"
$versionText = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5a3bac695026160bcf08146c67f957847f155c37/e2e/6d9ee57a-572e-4ce2-b42d-a8b49321472b.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/444133beedd7e169e1d11d60f530615a0b07b21d/e2e/6d9ee57a-572e-4ce2-b42d-a8b49321472b.md."

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("J7").Value = "6d9ee57a-572e-4ce2-b42d-a8b49321472b.fc3bfc6d36e4665b484beb8f95d6f18c90ef34c1.zh-cn.xlf"
$wsZh.Range("K7").Value = "2016-08-25 18:57:34"
$wsZh.Range("P7").Value = $versionText

$wsZh.Hyperlinks.Add($wsZh.Range("I7"), $latestUrl, "", "", "6d9ee57a-572e-4ce2-b42d-a8b49321472b.md")

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("J7").Value = "6d9ee57a-572e-4ce2-b42d-a8b49321472b.fc3bfc6d36e4665b484beb8f95d6f18c90ef34c1.de-de.xlf"
$wsDe.Range("K7").Value = "2016-08-25 18:57:41"
$wsDe.Range("P7").Value = $versionText

$wsDe.Hyperlinks.Add($wsDe.Range("I7"), $latestUrl, "", "", "6d9ee57a-572e-4ce2-b42d-a8b49321472b.md")
